$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.399.31"
$ws.Range("E2").Value = "  -0.94%  "
$ws.Range("D3").Value = "2.520.48"
$ws.Range("E3").Value = "  -0.77%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "542.76"
$ws.Range("E5").Value = "  -0.36%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.32"
$ws.Range("E6").Value = "  -3.65%  "
$ws.Range("E7").Value = "  +0.38%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.565"
$ws.Range("E8").Value = "  -1.26%  "
$ws.Range("D9").Value = "2.525.60"
$ws.Range("E9").Value = "  -1.72%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.101"
$ws.Range("E10").Value = "  +0.05%  "
$ws.Range("E11").Value = "  +0.38%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.42"
$ws.Range("E12").Value = "  -2.45%  "
$ws.Range("E13").Value = "  -3.29%  "
$ws.Range("D14").Value = "2.971.18"
$ws.Range("E14").Value = "  -0.59%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.39"
$ws.Range("E15").Value = "  -2.39%  "
$ws.Range("D16").Value = "59.280.84"
$ws.Range("E16").Value = "  -1.01%  "
$ws.Range("E17").Value = "  -1.50%  "
$ws.Range("D18").Value = "2.517.21"
$ws.Range("E18").Value = "  -2.18%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.11"
$ws.Range("E19").Value = "  -1.90%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.28"
$ws.Range("E20").Value = "  -1.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "325.59"
$ws.Range("E21").Value = "  -0.96%  "
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("E23").Value = "  -1.30%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.37"
$ws.Range("E24").Value = "  +1.60%  "
$ws.Range("E25").Value = "  -5.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.169"
$ws.Range("E26").Value = "  +0.87%  "
$ws.Range("E27").Value = "  +1.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.76"
$ws.Range("E28").Value = "  -3.46%  "
$ws.Range("D29").Value = "0.0₃0786"
$ws.Range("E29").Value = "  -1.63%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.72"
$ws.Range("E30").Value = "  -5.64%  "
$ws.Range("E31").Value = "  -1.37%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "164.98"
$ws.Range("E32").Value = "  +1.21%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.998"
$ws.Range("E33").Value = "  +0.13%  "
$ws.Range("E34").Value = "  -9.71%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.42"
$ws.Range("E35").Value = "  -3.81%  "
$ws.Range("E36").Value = "  -1.42%  "
$ws.Range("E37").Value = "  -6.29%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.59"
$ws.Range("E38").Value = "  -2.73%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.68"
$ws.Range("E39").Value = "  -1.86%  "
$ws.Range("E40").Value = "  -3.87%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.22"
$ws.Range("E41").Value = "  -8.64%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "281.66"
$ws.Range("E42").Value = "  -7.66%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.997"
$ws.Range("E43").Value = "  +0.44%  "
$ws.Range("E44").Value = "  -1.51%  "
$ws.Range("E45").Value = "  +0.31%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0514"
$ws.Range("E48").Value = "  -1.57%  "
$ws.Range("E49").Value = "  -2.53%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.94"
$ws.Range("E50").Value = "  -2.68%  "
$ws.Range("D51").Value = "1.772.81"
$ws.Range("E51").Value = "  -2.88%  "

# Row 46/47: Stellar and Aave swap places with updated values
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "125.72"
$ws.Range("E46").Value = "  +0.00%  "

$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0940"
$ws.Range("E47").Value = "  +0.06%  "
